$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: the "has_been_verified_by_agent" question now points at the
#     new "true_false_unsure" choice list instead of the old "true_false_01" one.
$survey.Range("E4").Value = "true_false_unsure"

# --- choices sheet: rename the "true_false_01" choice list to "true_false_unsure"
#     (same False/True rows) and add a new "Unsure" (-1) option, followed by a
#     blank spacer row before the male_female list.
$choices.Range("A6").Value = "true_false_unsure"
$choices.Range("A7").Value = "true_false_unsure"

$choices.Range("A8").Value = "true_false_unsure"
$choices.Range("B8").Value = -1
$choices.Range("B8").NumberFormat = "0"
$choices.Range("C8").Value = "Unsure"
$choices.Range("C8").NumberFormat = "@"

# Insert a new blank spacer row below the newly-added choice, matching the
# height of the row above it, pushing the male_female choices down by one row.
$choices.Rows.Item(9).Insert()
$choices.Rows.Item(9).RowHeight = 19
$choices.Range("B9").NumberFormat = "0"
$choices.Range("C9").NumberFormat = "@"

# --- view state: "survey" becomes the active/selected tab (was "choices"),
#     with a new selection on each sheet. Select on "choices" first so that
#     the final active sheet/selection ends up on "survey".
$choices.Activate()
$choices.Range("A9:XFD9").Select()

$survey.Activate()
$survey.Range("E5").Select()
